# Applies the "Schedule source completed" edit:
#  - Rename the four person labels (javi/pedro/cris/fran -> A/B/C/D)
#  - Update the numeric data in columns A and C
#  - Move the active cell/selection to G7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the D-column labels (shared strings javi/pedro/cris/fran -> A/B/C/D)
$ws.Range("D1").Value = "A"
$ws.Range("D2").Value = "B"
$ws.Range("D3").Value = "C"
$ws.Range("D4").Value = "D"

# Update column A values
$ws.Range("A2").Value = 8
$ws.Range("A3").Value = 12
$ws.Range("A4").Value = 80

# Update column C values
$ws.Range("C1").Value = 2
$ws.Range("C2").Value = 50
$ws.Range("C3").Value = 5

# Move the selection/active cell to G7
$ws.Range("G7").Select()
